$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add a new column BB (the next QoQ forecast vintage) and a new final
# row (83) to the combined QoQ component forecast sheet.
# ------------------------------------------------------------------

# Row 1 header: new vintage date in BB1. Give it the same date
# number-format / font / border as the rest of the header row by
# copying BA1's formatting onto it (value is set first so the copy
# only needs to bring across formatting, not clobber the new value).
$ws.Range("BB1").Value = 45986
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null

# Rows 2-82: column BB holds this vintage's forecast. For the older
# rows that's identical to the prior vintage (column BA); for the
# most recent rows the forecast itself has been revised.
$bbValues = @(0.4011221634494007, 0.4238702649943349, 1.065831052964711, 0.1679892790148898, 0.7826320556679605, 1.608891637697312, 0.7314738238097078, 0.6066829402737994, 0.3745978975504585, 1.488552458234295, -1.412768783584468, 1.150059426400716, 0.08121742761659334, -0.02717560123848273, 0.7243763014332956, 0.06292774364702325, 0.6288730841304044, -0.04932303278563666, 0.7008158264609534, -0.05880790937801805, 0.4315408898675201, 0.01220354597239748, -0.0390655612641666, 0.5861016707528819, -0.09710917920445183, 0.5608392729090639, 0.4327290387686133, 0.6415186033097626, 0.3424960270584307, 0.450007896627298, 0.5281644231130258, 0.9569320119344127, 1.458965617195048, 1.616667503385756, 0.7391298723574948, 0.05368498059650051, 0.2682953781150843, 0.07985748687143257, 0.5497494968945205, 0.5291048039492949, 0.4298198366907684, -0.00001256079109168695, 0.5955762646739942, 0.1202567733712812, 0.378819647862997, 1.1, 0.5, 1.4, 0.1, 1.14595180410879, 1.264351713671076, 0.9847891746257034, 0.1306165969516826, -0.702620495601451, 4.902896932015494, -3.137970849842105, 0.1542746451113572, 0.833545251617636, 0.3586840777896043, -1.11464790305466, 0.1937043270045251, -1.254351945331166, 0.05976999529040938, 1.365536114841987, 0.5893749959732304, 0.450398338429352, 1.478730722800918, 1.16573263774626, 0.4367731001324842, 0.1765865160815849, 0.2412052862208469, 0.768168485846715, 0.4322849306135402, 0.5070298915095297, 0.507551532258898, 0.5015060205255966, 0.5035546723560925, 0.503320960716943, 0.5032409843702599, 0.5032869668510612, 0.503277349074597)
for ($i = 0; $i -lt $bbValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 54).Value = $bbValues[$i]
}

# Row 83: a brand-new observation date, with the usual date formatting
# copied from the row above (A82), plus its BB forecast value.
$ws.Range("A83").Value = 46934
$ws.Range("A82").Copy() | Out-Null
$ws.Range("A83").PasteSpecial(-4122) | Out-Null
$ws.Range("BB83").Value = 0.503277065156336

"done"
